# Update loading_percent results for the 380 kV case (Case_0_167).
# Columns C:G and I:O on rows 2-25 are refreshed with new load-flow
# results; column H stays 0 and is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 9.428760145589459
$ws.Range("D2").Value = 8.050916247272447
$ws.Range("E2").Value = 13.13308168782856
$ws.Range("F2").Value = 37.77589305663253
$ws.Range("G2").Value = 3.690168876384484
$ws.Range("I2").Value = 29.35703054150622
$ws.Range("J2").Value = 10.28564733870244
$ws.Range("K2").Value = 16.03830506086214
$ws.Range("L2").Value = 10.26208512969714
$ws.Range("M2").Value = 18.59136596431139
$ws.Range("N2").Value = 19.14427219226823
$ws.Range("O2").Value = 29.09722937822778

# Row 3
$ws.Range("C3").Value = 9.40842994453245
$ws.Range("D3").Value = 8.03812260440194
$ws.Range("E3").Value = 13.15247976500198
$ws.Range("F3").Value = 37.86086263577771
$ws.Range("G3").Value = 3.692245676917338
$ws.Range("I3").Value = 29.43241442933885
$ws.Range("J3").Value = 10.30659931843412
$ws.Range("K3").Value = 15.72881507287311
$ws.Range("L3").Value = 10.28117286066747
$ws.Range("M3").Value = 18.4804777005878
$ws.Range("N3").Value = 19.18960923536597
$ws.Range("O3").Value = 29.184380304146

# Row 4
$ws.Range("C4").Value = 9.397338433290111
$ws.Range("D4").Value = 8.031172973431556
$ws.Range("E4").Value = 13.1658893008232
$ws.Range("F4").Value = 37.92109740921168
$ws.Range("G4").Value = 3.693589302120272
$ws.Range("I4").Value = 29.48445047007194
$ws.Range("J4").Value = 10.32024320410143
$ws.Range("K4").Value = 15.53756120404485
$ws.Range("L4").Value = 10.2936499158931
$ws.Range("M4").Value = 18.41435656598182
$ws.Range("N4").Value = 19.2192456570736
$ws.Range("O4").Value = 29.24367802000843

# Row 5
$ws.Range("C5").Value = 9.39317174400048
$ws.Range("D5").Value = 8.02857071466199
$ws.Range("E5").Value = 13.17173122963092
$ws.Range("F5").Value = 37.94766756448644
$ws.Range("G5").Value = 3.694154109348557
$ws.Range("I5").Value = 29.50709939961597
$ws.Range("J5").Value = 10.32599962621051
$ws.Range("K5").Value = 15.45941473875419
$ws.Range("L5").Value = 10.29892524600794
$ws.Range("M5").Value = 18.38792575092342
$ws.Range("N5").Value = 19.23177606711027
$ws.Range("O5").Value = 29.26929503656954

# Row 6
$ws.Range("C6").Value = 9.392501292723663
$ws.Range("D6").Value = 8.028152545995059
$ws.Range("E6").Value = 13.17272408739379
$ws.Range("F6").Value = 37.95220165214212
$ws.Range("G6").Value = 3.694248939752628
$ws.Range("I6").Value = 29.51094737661976
$ws.Range("J6").Value = 10.32696735334052
$ws.Range("K6").Value = 15.44642890709723
$ws.Range("L6").Value = 10.29981274931481
$ws.Range("M6").Value = 18.38356855534845
$ws.Range("N6").Value = 19.23388413816534
$ws.Range("O6").Value = 29.27363639102081

# Row 7
$ws.Range("C7").Value = 9.397280805578557
$ws.Range("D7").Value = 8.031136945530653
$ws.Range("E7").Value = 13.16596655826748
$ws.Range("F7").Value = 37.92144755350962
$ws.Range("G7").Value = 3.693596849314779
$ws.Range("I7").Value = 29.48475007835409
$ws.Range("J7").Value = 10.32032004126771
$ws.Range("K7").Value = 15.53650800656104
$ws.Range("L7").Value = 10.29372028758551
$ws.Range("M7").Value = 18.41399800332038
$ws.Range("N7").Value = 19.21941280988619
$ws.Range("O7").Value = 29.24401762132508

# Row 8
$ws.Range("C8").Value = 9.421463599229698
$ws.Range("D8").Value = 8.046318264085253
$ws.Range("E8").Value = 13.1394593192227
$ws.Range("F8").Value = 37.80351505009468
$ws.Range("G8").Value = 3.690870780564949
$ws.Range("I8").Value = 29.38182828904039
$ws.Range("J8").Value = 10.29271012370346
$ws.Range("K8").Value = 15.93190604565789
$ws.Range("L8").Value = 10.26850970979967
$ws.Range("M8").Value = 18.55273614507455
$ws.Range("N8").Value = 19.1595314881713
$ws.Range("O8").Value = 29.12607658048208

# Row 9
$ws.Range("C9").Value = 9.479768838080984
$ws.Range("D9").Value = 8.083182342915933
$ws.Range("E9").Value = 13.09935104444817
$ws.Range("F9").Value = 37.63637166905322
$ws.Range("G9").Value = 3.686065719127724
$ws.Range("I9").Value = 29.22571026648086
$ws.Range("J9").Value = 10.24472914490114
$ws.Range("K9").Value = 16.69334100413541
$ws.Range("L9").Value = 10.22505905669991
$ws.Range("M9").Value = 18.83948094342303
$ws.Range("N9").Value = 19.05634174672618
$ws.Range("O9").Value = 28.94081304345809

# Row 10
$ws.Range("C10").Value = 9.529015861786357
$ws.Range("D10").Value = 8.114460390939977
$ws.Range("E10").Value = 13.07709121387987
$ws.Range("F10").Value = 37.55283939951062
$ws.Range("G10").Value = 3.682861657355556
$ws.Range("I10").Value = 29.13898873133651
$ws.Range("J10").Value = 10.21320446187641
$ws.Range("K10").Value = 17.23887689355851
$ws.Range("L10").Value = 10.19675804998843
$ws.Range("M10").Value = 19.05784674888521
$ws.Range("N10").Value = 18.98915411957376
$ws.Range("O10").Value = 28.83289048086489

# Row 11
$ws.Range("C11").Value = 9.552758738476678
$ws.Range("D11").Value = 8.129570130267934
$ws.Range("E11").Value = 13.06852285830819
$ws.Range("F11").Value = 37.52339444267807
$ws.Range("G11").Value = 3.681474150340822
$ws.Range("I11").Value = 29.10563232008402
$ws.Range("J11").Value = 10.19966609117549
$ws.Range("K11").Value = 17.48298018804303
$ws.Range("L11").Value = 10.18466384745318
$ws.Range("M11").Value = 19.15858091573689
$ws.Range("N11").Value = 18.96045068309957
$ws.Range("O11").Value = 28.78994183879588

# Row 12
$ws.Range("C12").Value = 9.561937405008855
$ws.Range("D12").Value = 8.135415748696524
$ws.Range("E12").Value = 13.06550162157406
$ws.Range("F12").Value = 37.51347608506321
$ws.Range("G12").Value = 3.680958753321837
$ws.Range("I12").Value = 29.09387868621651
$ws.Range("J12").Value = 10.19465437075633
$ws.Range("K12").Value = 17.57474935586428
$ws.Range("L12").Value = 10.18019581320464
$ws.Range("M12").Value = 19.19690352407602
$ws.Range("N12").Value = 18.94984817117607
$ws.Range("O12").Value = 28.77456406947494

# Row 13
$ws.Range("C13").Value = 9.559952345801642
$ws.Range("D13").Value = 8.134151326411331
$ws.Range("E13").Value = 13.06614237230879
$ws.Range("F13").Value = 37.51555737386837
$ws.Range("G13").Value = 3.681069308353344
$ws.Range("I13").Value = 29.09637098415217
$ws.Range("J13").Value = 10.19572862776158
$ws.Range("K13").Value = 17.5550161785813
$ws.Range("L13").Value = 10.18115311922602
$ws.Range("M13").Value = 19.18864259269158
$ws.Range("N13").Value = 18.95211975406623
$ws.Range("O13").Value = 28.77783651135326

# Row 14
$ws.Range("C14").Value = 9.553510139199304
$ws.Range("D14").Value = 8.130048586498248
$ws.Range("E14").Value = 13.06826982529162
$ws.Range("F14").Value = 37.52255375811663
$ws.Range("G14").Value = 3.681431547747945
$ws.Range("I14").Value = 29.10464774328961
$ws.Range("J14").Value = 10.19925147216133
$ws.Range("K14").Value = 17.49054392762812
$ws.Range("L14").Value = 10.18429402164218
$ws.Range("M14").Value = 19.16173033728452
$ws.Range("N14").Value = 18.95957306373663
$ws.Range("O14").Value = 28.78865893278256

# Row 15
$ws.Range("C15").Value = 9.549588401890336
$ws.Range("D15").Value = 8.127551587928535
$ws.Range("E15").Value = 13.06960202868035
$ws.Range("F15").Value = 37.52699970562547
$ws.Range("G15").Value = 3.681654733507336
$ws.Range("I15").Value = 29.10983184106286
$ws.Range("J15").Value = 10.20142427597885
$ws.Range("K15").Value = 17.45096350947285
$ws.Range("L15").Value = 10.18623246017309
$ws.Range("M15").Value = 19.14526809795914
$ws.Range("N15").Value = 18.9641731636066
$ws.Range("O15").Value = 28.79540341719779

# Row 16
$ws.Range("C16").Value = 9.527490757549291
$ws.Range("D16").Value = 8.113490426771115
$ws.Range("E16").Value = 13.0776824723776
$ws.Range("F16").Value = 37.55493598854426
$ws.Range("G16").Value = 3.682953739286488
$ws.Range("I16").Value = 29.14129140058022
$ws.Range("J16").Value = 10.21410533652454
$ws.Range("K16").Value = 17.22283518008958
$ws.Range("L16").Value = 10.19756409766639
$ws.Range("M16").Value = 19.05128960025796
$ws.Range("N16").Value = 18.99106733641287
$ws.Range("O16").Value = 28.83582113621643

# Row 17
$ws.Range("C17").Value = 9.514274357967494
$ws.Range("D17").Value = 8.105088079861563
$ws.Range("E17").Value = 13.08303813140269
$ws.Range("F17").Value = 37.57426632729158
$ws.Range("G17").Value = 3.683768540054965
$ws.Range("I17").Value = 29.16215265725469
$ws.Range("J17").Value = 10.22208996657176
$ws.Range("K17").Value = 17.08178432844786
$ws.Range("L17").Value = 10.20471520195574
$ws.Range("M17").Value = 18.99397778219078
$ws.Range("N17").Value = 19.00804207630219
$ws.Range("O17").Value = 28.86219180339649

# Row 18
$ws.Range("C18").Value = 9.506799087763552
$ws.Range("D18").Value = 8.100338391095557
$ws.Range("E18").Value = 13.08626522515378
$ws.Range("F18").Value = 37.5861897258084
$ws.Range("G18").Value = 3.684243787026078
$ws.Range("I18").Value = 29.17472499022371
$ws.Range("J18").Value = 10.22675806140923
$ws.Range("K18").Value = 17.00027777183447
$ws.Range("L18").Value = 10.20890177500583
$ws.Range("M18").Value = 18.9611466443605
$ws.Range("N18").Value = 19.01798066161193
$ws.Range("O18").Value = 28.87793787144837

# Row 19
$ws.Range("C19").Value = 9.504289955625216
$ws.Range("D19").Value = 8.098744588156796
$ws.Range("E19").Value = 13.0873830701418
$ws.Range("F19").Value = 37.59036499996932
$ws.Range("G19").Value = 3.684405831715309
$ws.Range("I19").Value = 29.17908021563084
$ws.Range("J19").Value = 10.22835158599784
$ws.Range("K19").Value = 16.97261877155027
$ws.Range("L19").Value = 10.2103319033969
$ws.Range("M19").Value = 18.95005418343669
$ws.Range("N19").Value = 19.02137580170383
$ws.Range("O19").Value = 28.88336848411124

# Row 20
$ws.Range("C20").Value = 9.515668214051326
$ws.Range("D20").Value = 8.105973941668667
$ws.Range("E20").Value = 13.08245283773602
$ws.Range("F20").Value = 37.57212524101833
$ws.Range("G20").Value = 3.683681120911161
$ws.Range("I20").Value = 29.15987257555099
$ws.Range("J20").Value = 10.22123217326215
$ws.Range("K20").Value = 17.09683915485343
$ws.Range("L20").Value = 10.20394635608039
$ws.Range("M20").Value = 19.00006512280567
$ws.Range("N20").Value = 19.00621696249919
$ws.Range("O20").Value = 28.85932472669201

# Row 21
$ws.Range("C21").Value = 9.555397316444312
$ws.Range("D21").Value = 8.131250321454303
$ws.Range("E21").Value = 13.06763888225566
$ws.Range("F21").Value = 37.52046530836519
$ws.Range("G21").Value = 3.681324877660965
$ws.Range("I21").Value = 29.10219282706841
$ws.Range("J21").Value = 10.19821361046578
$ws.Range("K21").Value = 17.50949974245195
$ws.Range("L21").Value = 10.18336843197802
$ws.Range("M21").Value = 19.16963051667991
$ws.Range("N21").Value = 18.95737660899849
$ws.Range("O21").Value = 28.78545606425724

# Row 22
$ws.Range("C22").Value = 9.582454505836807
$ws.Range("D22").Value = 8.148490623029184
$ws.Range("E22").Value = 13.0592590352166
$ws.Range("F22").Value = 37.49388291840629
$ws.Range("G22").Value = 3.679843328410255
$ws.Range("I22").Value = 29.06961228865675
$ws.Range("J22").Value = 10.18383957510757
$ws.Range("K22").Value = 17.77527187536176
$ws.Range("L22").Value = 10.17057091485114
$ws.Range("M22").Value = 19.28147020872019
$ws.Range("N22").Value = 18.92701175271468
$ws.Range("O22").Value = 28.74234351354234

# Row 23
$ws.Range("C23").Value = 9.567915357671557
$ws.Range("D23").Value = 8.139224165373269
$ws.Range("E23").Value = 13.06361259362962
$ws.Range("F23").Value = 37.50741300036495
$ws.Range("G23").Value = 3.680628732488499
$ws.Range("I23").Value = 29.08653255120567
$ws.Range("J23").Value = 10.19145010508271
$ws.Range("K23").Value = 17.63380961727169
$ws.Range("L23").Value = 10.17734172190404
$ws.Range("M23").Value = 19.22169412046537
$ws.Range("N23").Value = 18.94307597529304
$ws.Range("O23").Value = 28.7648802258343

# Row 24
$ws.Range("C24").Value = 9.515037668181566
$ws.Range("D24").Value = 8.105573191249453
$ws.Range("E24").Value = 13.08271698779141
$ws.Range("F24").Value = 37.57309070273981
$ws.Range("G24").Value = 3.683720621898897
$ws.Range("I24").Value = 29.16090159729213
$ws.Range("J24").Value = 10.22161973979061
$ws.Range("K24").Value = 17.09003415265477
$ws.Range("L24").Value = 10.20429371669329
$ws.Range("M24").Value = 18.99731266610171
$ws.Range("N24").Value = 19.00704153700274
$ws.Range("O24").Value = 28.86061910928878

# Row 25
$ws.Range("C25").Value = 9.46285280577453
$ws.Range("D25").Value = 8.072463201127295
$ws.Range("E25").Value = 13.10893336433004
$ws.Range("F25").Value = 37.6747036425048
$ws.Range("G25").Value = 3.687308081468866
$ws.Range("I25").Value = 29.26303835072717
$ws.Range("J25").Value = 10.25705266203884
$ws.Range("K25").Value = 16.48941734016732
$ws.Range("L25").Value = 10.2361755230635
$ws.Range("M25").Value = 18.76046478401508
$ws.Range("N25").Value = 19.08273871300532
$ws.Range("O25").Value = 28.98599054945704
